$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9

# Row 6
$ws.Range("G6").Value = 1.7
$ws.Range("L6").Value = 5.5
$ws.Range("M6").Value = 1.08
$ws.Range("N6").Value = 7.5
$ws.Range("Q6").Value = 2.3
$ws.Range("R6").Value = 1.6
$ws.Range("AV6").Value = 67

# Row 9
$ws.Range("S9").Value = 1.4

# Row 10
$ws.Range("S10").Value = 1.33
